$d = $word.ActiveDocument

# Remove the old "_GoBack" bookmark first (it currently sits in the very
# last paragraph of the document). It will be re-created below at the
# new edit location, mirroring what Word itself does when you type in a
# fresh spot: the lone "_GoBack" bookmark follows the latest edit.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Locate the empty "ListParagraph" paragraph that immediately precedes
# the "Replace the paths referenced..." list item (numId=3). Found by
# searching for that list item's text and walking back one paragraph,
# which is more robust than depending on a hard-coded paragraph index.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "Replace the paths referenced in the command*") {
        $target = $d.Paragraphs.Item($i - 1)
        break
    }
}

$curlyOpen  = [char]0x201C
$curlyClose = [char]0x201D

# Fill the previously-empty paragraph with the new sentence, splitting
# out the word "src" into its own run (bracketed by spell-check
# proof-error markers, matching how Word flags a word not in its
# dictionary), and drop the "_GoBack" bookmark at the end of it (the
# newest edit location).
$r = $target.Range
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr>' +
    '<w:r><w:t>Usually run the instrument the java source code folder is fine, for example, the ' + $curlyOpen + '</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>src</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>' + $curlyClose + ' folder.</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '</w:p>'
[void]$r.InsertXML($newParaXml)

# Insert a new, empty "ListParagraph" paragraph right after it, so the
# numbered item "Replace the paths referenced..." keeps a blank spacer
# paragraph in front of it, just like the rest of the list.
$endOfFilled = $d.Range($target.Range.End, $target.Range.End)
$blankParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr>' +
    '</w:p>'
[void]$endOfFilled.InsertXML($blankParaXml)
